$d = $word.ActiveDocument

$replacements = @(
    @("927×5=", "294×9="),
    @("942×9=", "446×4="),
    @("178×6=", "805×5="),
    @("949×7=", "739×2="),
    @("146×8=", "736×3="),
    @("648×6=", "631×4="),
    @("114×9=", "809×4="),
    @("696×7=", "578×4="),
    @("817×4=", "831×3="),
    @("518×7=", "335×6="),
    @("406×7=", "993×5="),
    @("221×6=", "464×9="),
    @("818×4=", "947×4="),
    @("849×9=", "855×6="),
    @("111×2=", "724×2="),
    @("545×9=", "404×4="),
    @("265×7=", "525×3="),
    @("975×6=", "865×5="),
    @("255×2=", "283×4="),
    @("734×5=", "900×8="),
    @("811×8=", "153×2="),
    @("400×5=", "287×8="),
    @("905×4=", "435×5="),
    @("589×5=", "364×4="),
    @("566×7=", "691×6=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
